$d = $word.ActiveDocument

# ----------------------------------------------------------------------------
# Task 5 ("Team Meetings") hours: "6h" -> "4h".
# The diff shows the original single run "6h" split into two runs "4" + "h"
# (both keeping the existing single-underline formatting), so we replace the
# two characters individually rather than doing a whole-string replace (which
# would collapse back into one run).
# ----------------------------------------------------------------------------
$rng1 = $d.Content
$found1 = $rng1.Find.Execute("Hours:_____6h", $true, $false, $false, $false, $false,
                              $true, 1, $false, "", 0)
if (-not $found1) {
    throw "Could not find 'Hours:_____6h' (Task 5 hours) in the document"
}
$hoursStart = $rng1.End - 2   # position of the "6" in "6h"

$r1 = $d.Range($hoursStart, $hoursStart + 1)
$r1.Bold = 1
$r1.Text = "4"
$r1.Bold = 0

$r2 = $d.Range($hoursStart + 1, $hoursStart + 2)
$r2.Bold = 1
$r2.Text = "h"
$r2.Bold = 0

# ----------------------------------------------------------------------------
# Total Hours: "23h*" -> "21h*".
# Only the "3h*" run is touched by the diff; it becomes two runs "1" + "h*"
# (the preceding "2" run stays as-is, still its own separate run).
# ----------------------------------------------------------------------------
$rng2 = $d.Content
$found2 = $rng2.Find.Execute("Total Hours:  _____23h*", $true, $false, $false, $false, $false,
                              $true, 1, $false, "", 0)
if (-not $found2) {
    throw "Could not find 'Total Hours:  _____23h*' in the document"
}
$totalStart = $rng2.End - 3   # position of the "3" in "3h*"

# Re-write the preceding "2" run in place so it stays a run distinct from the
# new "1" run that follows it (both share the same underline formatting, and
# would otherwise be recombined into a single "21" run).
$rPrev = $d.Range($totalStart - 1, $totalStart)
$rPrev.Bold = 1
$rPrev.Text = "2"
$rPrev.Bold = 0

$r3 = $d.Range($totalStart, $totalStart + 1)
$r3.Bold = 1
$r3.Text = "1"
$r3.Bold = 0

$r4 = $d.Range($totalStart + 1, $totalStart + 3)
$r4.Bold = 1
$r4.Text = "h*"
$r4.Bold = 0
